# Update scripts with new TPM values for Ereg-Egfr (FAPs -> *) pairs.
# The old "MuSCs" sending-cluster block (rows 6-9) is dropped entirely,
# and the remaining FAPs-sending rows (2-5) get refreshed TPM-derived
# statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "MuSCs" sending-cluster rows (old rows 6-9).
$ws.Rows.Item(6).Resize(4).Delete()

# Row 2: FAPs -> Ereg/Egfr -> ECs
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 0.1120545175691111
$ws.Range("R2").Value = 1.008490658122
$ws.Range("S2").Value = 0.01103063309339269
$ws.Range("T2").Value = 0.01103063309339269

# Row 3: FAPs -> Ereg/Egfr -> FAPs
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("S3").Value = 0.7476219244149905
$ws.Range("T3").Value = 0.7476219244149904

# Row 4: FAPs -> Ereg/Egfr -> MuSCs
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 2.429169069837333
$ws.Range("R4").Value = 21.862521628536
$ws.Range("S4").Value = 0.2391271080585153
$ws.Range("T4").Value = 0.2391271080585153

# Row 5: FAPs -> Ereg/Egfr -> Resolving-Mac
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.02255523337933333
$ws.Range("R5").Value = 0.202997100414
$ws.Range("S5").Value = 0.002220334433101459
$ws.Range("T5").Value = 0.002220334433101458
